# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value for sheet "展览"
$updates1 = @{
    6  = 670
    7  = 270
    12 = 3373
    13 = 105
    14 = 76
    18 = 572
    20 = 673
    21 = 198
    22 = 110
    24 = 49
    26 = 2403
    27 = 4942
    31 = 1268
    32 = 270
    33 = 2191
    37 = 76
    41 = 771
    42 = 23
    43 = 448
    45 = 454
}

# Map of row -> new F-column value for sheet "全部类型"
$updates4 = @{
    6  = 670
    7  = 270
    12 = 3373
    13 = 105
    14 = 76
    19 = 572
    21 = 673
    22 = 198
    23 = 110
    25 = 49
    27 = 2403
    28 = 4942
    32 = 1268
    33 = 270
    34 = 2191
    38 = 76
    42 = 771
    43 = 23
    44 = 448
    46 = 454
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
